# Add a new paragraph after the final "Andy Herrold" paragraph containing
# the game description text (with spell-check proofErr markers around the
# two coined/proper-noun-like words "pygame" and "simpleGE").

$d = $word.ActiveDocument

# Build the new paragraph's WordprocessingML. Every run uses the same
# Times New Roman rFonts-only formatting as the rest of the document, and
# no paragraph-level alignment/bold/size is set (plain body paragraph).
$newParaXml = @'
<w:p>
  <w:pPr>
    <w:rPr>
      <w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/>
    </w:rPr>
  </w:pPr>
  <w:r>
    <w:rPr>
      <w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/>
    </w:rPr>
    <w:t xml:space="preserve">“Beat the Dealer” will be a simple 2D card game utilizing </w:t>
  </w:r>
  <w:proofErr w:type="spellStart"/>
  <w:r>
    <w:rPr>
      <w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/>
    </w:rPr>
    <w:t>pygame</w:t>
  </w:r>
  <w:proofErr w:type="spellEnd"/>
  <w:r>
    <w:rPr>
      <w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/>
    </w:rPr>
    <w:t xml:space="preserve"> and </w:t>
  </w:r>
  <w:proofErr w:type="spellStart"/>
  <w:r>
    <w:rPr>
      <w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/>
    </w:rPr>
    <w:t>simpleGE</w:t>
  </w:r>
  <w:proofErr w:type="spellEnd"/>
  <w:r>
    <w:rPr>
      <w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/>
    </w:rPr>
    <w:t xml:space="preserve">. It is based upon the popular casino game blackjack. The player will be represented by two random cards at the bottom of the screen. The dealer will be represented by two random cards near the top of the screen. The players objective will be to make a hand with a higher total than that of the dealer without exceeding a total of 21 (“busting”). A casino or blackjack table image will make up the background. </w:t>
  </w:r>
</w:p>
'@

# Wrap it in the minimal single-part package InsertXML expects.
$pkg = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?>' +
       '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">' +
       '<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">' +
       '<pkg:xmlData>' +
       '<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">' +
       '<w:body>' + $newParaXml + '</w:body>' +
       '</w:document>' +
       '</pkg:xmlData></pkg:part></pkg:package>'

# Insert it as a brand-new paragraph at the very end of the document, after
# the existing (untouched) trailing "Andy Herrold" paragraph.
$endRange = $d.Content
$endRange.Collapse(0)
[void]$endRange.InsertXML($pkg)
